# Update Name of Algo
# Apply the recomputed values for the RandomForest result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -13.78619999999998
$ws.Range("E4").Value = 13.9104

$ws.Range("E5").Value = 12.87919999999999

$ws.Range("C6").Value = -11.66780000000001
$ws.Range("E6").Value = 12.39280000000001

$ws.Range("C7").Value = -12.0666

$ws.Range("C8").Value = -11.9721
$ws.Range("E8").Value = 13.85389999999999

$ws.Range("C16").Value = -11.79860000000001
$ws.Range("E16").Value = 12.6506

$ws.Range("C20").Value = -14.3589

$ws.Range("C21").Value = -13.09030000000001

$ws.Range("E22").Value = 12.13049999999999
